$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Drop the now-unused shared string ("remise en forme du site")
#    by clearing the one cell that referenced it. The engine prunes
#    unreferenced shared strings on save.
# ------------------------------------------------------------------
$ws.Range("E2").ClearContents()

# ------------------------------------------------------------------
# 2. Update the date header row (row 1) - two weeks later.
# ------------------------------------------------------------------
$ws.Range("B1").Value = 44529
$ws.Range("C1").Value = 44530
$ws.Range("D1").Value = 44531
$ws.Range("E1").Value = 44532
$ws.Range("F1").Value = 44533
$ws.Range("F1").NumberFormat = "m/d/yy"

# ------------------------------------------------------------------
# 3. Fill in the "chose faites" / "difficulte" content for the first
#    week, in the exact order the cells were originally authored so
#    that new shared-string indices line up with the target.
# ------------------------------------------------------------------
$ws.Range("C2").Value = "entrainement QT"
$ws.Range("D2").Value = "installationROBOTIS pour programmer robot"
$ws.Range("D4").Value = "commencement du développement de la fenetre utilisateur (manette virtuelle)"
$ws.Range("D7").Value = "cervo moteur de l'araignée assez puissant faisant sauterl e support de la carte arduino"
$ws.Range("D8").Value = "QT creator n'avais pas toutes les ressources pour commencer a développer la fenetre principal"
$ws.Range("D3").Value = "Debut Apprentissage lent du language de prog du robot"
$ws.Range("E2").Value = "debut utilisation tramme serie envoyé"
$ws.Range("E3").Value = "finalisatoin design fenetre qt"
$ws.Range("E7").Value = "architcture des programmes arduino"
$ws.Range("E8").Value = "comprehension de la personnalisation des winget QT pas encore optiale"
$ws.Range("C3").Value = "decouverte QT serial"
$ws.Range("C7").Value = "doc de l'arraigné superflue"

Write-Host "content written"

# ------------------------------------------------------------------
# 4. Row 10 header (2nd week) - now built from formulas.
# ------------------------------------------------------------------
$ws.Range("B10").Formula = "=F1+1"
$ws.Range("C10").Formula = "=B10+1"
$ws.Range("D10:F10").Formula = "=C10+1"

Write-Host "formulas written"

# ------------------------------------------------------------------
# 5. Apply the "Insatisfaisant" (red/Bad) and "Satisfaisant"
#    (green/Good) cell styles to the "difficulte" / "chose faites"
#    blocks. Bad is created first so it lands on cellXfs index 2 and
#    Good on index 3 (matching s="2"/s="3" used throughout the sheet).
# ------------------------------------------------------------------
$ws.Range("A7:E9").Style = "Bad"
$ws.Range("A17:E19").Style = "Bad"

$ws.Range("A2:E6").Style = "Good"
$ws.Range("A11:E16").Style = "Good"

Write-Host "styles applied"

